$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

# ---------------------------------------------------------------------
# Paragraph 1 (the title) originally reads:
#   "Initio Simulator Programming: WS17 Sample Answers and Trouble Shooting"
# made up of 7 runs:
#   r1 "Initio"          [0,6)
#   r2 " "                [6,7)
#   r3 "Simulator "      [7,17)
#   r4 "Programming: "   [17,30)
#   r5 "WS1"             [30,33)
#   r6 "7"               [33,34)
#   r7 " Sample Answers and Trouble Shooting" [34,71)
#
# r1..r4 share one rPr (no color); r5..r7 share another (with color).
# The engine merges adjacent runs that end up with identical effective
# rPr whenever a Range.Text mutation touches their paragraph, so editing
# r1..r3 in place would otherwise also fuse r4 into them and fuse
# r5/r6/r7 together. To avoid that, every other run in each same-rPr
# group is given a momentarily different font size (a no-op once it is
# set back to the original 20pt / w:sz 40) so no two adjacent runs ever
# look identical while the text edits happen.
# ---------------------------------------------------------------------

$mark1 = $d.Range(0, 6);   $mark1.Font.Size = 21   # r1 "Initio"
$mark2 = $d.Range(7, 17);  $mark2.Font.Size = 21   # r3 "Simulator "
$mark3 = $d.Range(30, 33); $mark3.Font.Size = 21   # r5 "WS1"

# r1: "Initio" -> "Virtual "
$r1 = $d.Range(0, 6)
$r1.Text = "Virtual "

# r2: " " -> "Initio"  (now located right after "Virtual ", at offset 8)
$r2 = $d.Range(8, 9)
$r2.Text = "Initio"

# r3: "Simulator " -> " "  (now located at offset 14, length 10)
$r3 = $d.Range(14, 24)
$r3.Text = " "

# Restore the original font size (20pt == w:sz 40) everywhere we marked.
$reset1 = $d.Range(0, 8);   $reset1.Font.Size = 20   # "Virtual "
$reset2 = $d.Range(8, 14);  $reset2.Font.Size = 20   # "Initio"
$reset3 = $d.Range(14, 15); $reset3.Font.Size = 20   # " "
$reset4 = $d.Range(28, 31); $reset4.Font.Size = 20   # "WS1"

Write-Host "Title now reads:[" $p1.Range.Text "]"

# ---------------------------------------------------------------------
# Move the _GoBack bookmark from the end of the document (right after
# "initio.irRight()") to right after the new third run (the lone space
# that used to be "Simulator "), i.e. right before "Programming: ".
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$bmRange = $d.Range(15, 15)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "Done."
